# edit.ps1 - transforms Review 129 (Watch Your Steps) into Review 128 (CyBERT)
# following the unified diff supplied for Review_129.docx -> Review_128 content.
$d = $word.ActiveDocument

# Paragraph 1 (Heading1): review title line + arxiv link line
$oldTitle = "Review 129: [Short] Watch Your Steps: Local Image and Scene Editing by Text Instructions,  26.08.23https://arxiv.org/abs/2211.09800.pdf"
$newTitle = "Review 128: [Short] 18.08.23: CyBERT: Contextualized Embeddings for the Cybersecurity Domain (סקירה זו נכתבה על ידי עדן יבין)^lhttps://mdsoar.org/bitstream/handle/11603/25498/1117.pdf"
$foundTitle = $d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)
if (-not $foundTitle) { throw "Replacement failed: Paragraph 1 (Heading1): review title line + arxiv link line" }

# Paragraph 2 (bold): Paper link
$oldPaper = "Paper: https://arxiv.org/abs/2308.08947v1"
$newPaper = "Paper: https://arxiv.org/abs/2210.08218v1"
$foundPaper = $d.Content.Find.Execute($oldPaper, $true, $false, $false, $false, $false, $true, 1, $false, $newPaper, 2)
if (-not $foundPaper) { throw "Replacement failed: Paragraph 2 (bold): Paper link" }

# Paragraph 5: replace the whole review body with the new CyBERT summary (2 parts, double line-break)
$oldP5 = "מודלי הדיפוזיה לגמרי השתלטו כמעט על כל המשימות של הראייה הממוחשבת. למשל עריכה של תמונות (למשל להחליף ציפור בפרפר)בהתאם לתיאור טקסטואלי כבר מזמן עושים רק באמצעות מודלי דיפוזיה חזקים כמו  InstructPix2Pix או IP2P בקצרה. למרות התוצאות המדהימות עדיין יש אי התאמות בין התמונה הערוכה לבין המקורית. היום ב-#shorthebrewpapereviews נסקור מאמר שמנסה לתקן את אי דיוקים אלו בצורה די אלגנטית. בשלב הראשון המודל המוצע מאתר את מיקום הפיקסלים שאותם צריך לשנות(מסכה) ובשלב השני עורכים את התמונה רק באזורים של המסכה. כל זה נעשה באמצעות מודלי דיפוזיה באופן די אלגנטי. בשלב הראשון מרעישים את התמונה המקורית (עד רמת רעש מסוימת שהיא מהווה הייפרפרמטר חשוב מאוד) משתמשים במודל IP2P כדי לשערך את הרעש נוסף עבור ללא תופסת טקסט לעריכה ויחד איתו. כלומר במקרה הראשון אנו מפעילים מודל דיפוזיה סטנדרטי (ללא עריכה) ובמקרה השני כן עורכים את התמונה בהתאם לתיאור הטקסטואלי. לאחר מכן מחשבים את הערך המוחלט של ההפרש בין השערוכים אלו, מקצצים את החריגים (עם IQR עם מקדם 1.5). המסכה מקבלת ערך 1 (פיקסלים לעריכה) במקומות שההפרש הזה עולה על סף מסוים (הייפרפרמטר נוסף). בשלב השני מרעישים את התמונה (רמת הרעש עוד הייפרפרמטר). ואז באמצעות מסירים את הרעש עם מודל IP2P (עם תיאור טקסטואלי) באיזורים של המסכה ובכל האזורים האחרים עושים זאת עם מודל דיפוזיה רגיל (הטקסט המוסף הוא ריק). בנוסף המחברים מכלילים את הגישה שלהם ל-NeRF (ייצוג של מודלי 3D). בגדול עושים את מה שמתואר למעלה על views מכל הזווית תוך שמירה של קוהרנטיות ביניהם."
$newP5 = "מודלי שפה הראו את היכולת שלהם לעבוד בהמון תחומים בתוך עיבוד שפה טבעית. אחד התחומים שבהם ציפו להשפעה גדולה של מודלים אלו הוא תחום אבטחת המידע או בשמו היותר מוקר תחום הסייבר. במאמר ששמו הינו CyBERT הראו לראשונה את השילוב של מודלי שפה גדולים (מודל BERT) בתחום הסייבר. ^l^lהחוקרים ניסו להראות איך שיפור של BERT הקיים יכול להביא לשיפור ניכר של אותו מודל במשימות שונות כגון זיהוי אובייקטים הקשורים לתחום הסייבר או סיווג של מילה לאובייקט המתאים לה. למשל, ניתן לראות תוצאות של המשימה האחרונה המוזכרת בה המודל נאלץ לסווג מילה מתחום הסייבר לאובייקט המתאים לה. איך ביצעו זאת? על ידי הרחבת האימון של BERT עם מספר שלבים נוספים: "
$foundP5 = $d.Content.Find.Execute($oldP5, $true, $false, $false, $false, $false, $true, 1, $false, $newP5, 2)
if (-not $foundP5) { throw "Replacement failed: Paragraph 5: replace the whole review body with the new CyBERT summary (2 parts, double line-break)" }

# Paragraph 6: was a single empty run; fill it in with the new bullet list content.
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "– אוספים מסמכים רבים מתחום הסייבר – מנקים את המסמכים והופכים אותם לרשימה של טוקנים– את הרשימה של טוקנים מוספים למילון של ה-Tokenizer של BERT. בנוסף, מוסיפים למטריצת ה-Embedding הרגילה של BERT את הטוקנים עם ערך רנדומלי. מבצעים אימון נוסף של Masked Language Modeling, בה המודל נדרש להשלים מילות חסרות במשפטים. ככה הערכים הרנדומלים מתעדכנים. – מבצעים אימון נוסף בהתאם למשימה, למשל אימון מפוקח של זיהוי אובייקטים של סייבר במשפט. המאמר מעניין אותי במיוחד לא רק בשל היותו בתחום הסייבר אלא בשל התזכורת שהוא נותן כיצד ניתן לבצע אימון נוסף של מודלי שפה גדולים על תחומים חדשים."

# re-fetch the paragraph (its Range identity can go stale after the text assignment)
# then split the flat text into 4 segments joined by manual line breaks, in one pass
# (this keeps everything inside a single <w:r> and lets xml:space be recalculated correctly).
$p6 = $d.Paragraphs(6)
$flatP6 = "– אוספים מסמכים רבים מתחום הסייבר – מנקים את המסמכים והופכים אותם לרשימה של טוקנים– את הרשימה של טוקנים מוספים למילון של ה-Tokenizer של BERT. בנוסף, מוסיפים למטריצת ה-Embedding הרגילה של BERT את הטוקנים עם ערך רנדומלי. מבצעים אימון נוסף של Masked Language Modeling, בה המודל נדרש להשלים מילות חסרות במשפטים. ככה הערכים הרנדומלים מתעדכנים. – מבצעים אימון נוסף בהתאם למשימה, למשל אימון מפוקח של זיהוי אובייקטים של סייבר במשפט. המאמר מעניין אותי במיוחד לא רק בשל היותו בתחום הסייבר אלא בשל התזכורת שהוא נותן כיצד ניתן לבצע אימון נוסף של מודלי שפה גדולים על תחומים חדשים."
$brokenP6 = "– אוספים מסמכים רבים מתחום הסייבר ^l– מנקים את המסמכים והופכים אותם לרשימה של טוקנים^l– את הרשימה של טוקנים מוספים למילון של ה-Tokenizer של BERT. בנוסף, מוסיפים למטריצת ה-Embedding הרגילה של BERT את הטוקנים עם ערך רנדומלי. מבצעים אימון נוסף של Masked Language Modeling, בה המודל נדרש להשלים מילות חסרות במשפטים. ככה הערכים הרנדומלים מתעדכנים.^l – מבצעים אימון נוסף בהתאם למשימה, למשל אימון מפוקח של זיהוי אובייקטים של סייבר במשפט. המאמר מעניין אותי במיוחד לא רק בשל היותו בתחום הסייבר אלא בשל התזכורת שהוא נותן כיצד ניתן לבצע אימון נוסף של מודלי שפה גדולים על תחומים חדשים."
$foundP6Breaks = $p6.Range.Find.Execute($flatP6, $true, $false, $false, $false, $false, $true, 1, $false, $brokenP6, 2)
if (-not $foundP6Breaks) { throw "Replacement failed: paragraph 6 line breaks" }

Write-Host "All replacements applied successfully."
